$wb = $excel.ActiveWorkbook

$wsReport = $wb.Worksheets.Item("Active Report")
$wsDataSources = $wb.Worksheets.Item("Active DataSources")

# Sheet1 "Active Report": remove the "report1" row (row 2), leaving only the
# "report2" row, which moves up to row 2.
$wsReport.Rows.Item(2).Delete()

# Sheet2 "Active DataSources": update the datasource name/description.
$wsDataSources.Range("A2").Value = "postgres"
$wsDataSources.Range("B2").Value = "Datasource di produzione schema MONET"
$wsDataSources.Columns.Item(2).EntireColumn.AutoFit()
